# Denmark Division 3 - swap mixed-up rows back into their correct order.
# Each pair below represents two fixtures that were entered on the wrong
# row (same kick-off date/time); the full data for each match (id column
# "A" excluded) needs to be swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(2, 3),
    @(5, 6),
    @(66, 67),
    @(116, 117),
    @(124, 125),
    @(128, 129),
    @(132, 133),
    @(134, 136),
    @(140, 141),
    @(177, 178),
    @(194, 195)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
